$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "42 - Deep Thought"
$ws.Range("A17").HorizontalAlignment = -4131
$ws.Range("A18").Select() | Out-Null
